$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix / update existing row 2 (Energisa) text fields ---
$ws.Range("F2").Value = "The Energisa Group was founded in 1905 in Cataguases (MG) as Companhia Força e Luz Cataguazes-Leopoldina.It is one of the oldest groups in the Brazilian electricity sector.It grew through acquisitions and today operates in Distribution, Transmission and Generation in 11 states, focusing on solutions and diversification (Voltz, natural gas, solar)."
$ws.Range("I2").Value = "Big Company"
$ws.Range("J2").Value = "Cataguases benefits from the Energisa (CSE) headquarters with jobs, income and booming commerce.The group also invests in local education."
$ws.Range("K2").Value = "ABRADEE Award (quality and management) and recognized as GPTW (Best Company to Work For)."
$ws.Range("L2").Value = "Domestic (national) demand."

# --- Append new company rows 6-14 ---
$ws.Range("A6").Value = "Miba "
$ws.Range("B6").Value = "Av. Manoel Inácio Peixoto, 2147 - Industrial,"
$ws.Range("C6").Value = -21.4071436278473
$ws.Range("D6").Value = -42.6745858186452

$ws.Range("A7").Value = "Mais Energia Solar"
$ws.Range("B7").Value = "R. Ten. Fortunato, 183 - Centro"
$ws.Range("C7").Value = -21.3909432057039
$ws.Range("D7").Value = -42.6930428478218

$ws.Range("A8").Value = "Industrial Cataguases"
$ws.Range("B8").Value = "Praça José Ignácio Peixoto, 28 - Vila Tereza"
$ws.Range("C8").Value = -21.392621051068
$ws.Range("D8").Value = -42.6938404230373

$ws.Range("A9").Value = "Hidroazul"
$ws.Range("B9").Value = "R. João Dias Neto, 18 - Vila Reis"
$ws.Range("C9").Value = -21.3862116993511
$ws.Range("D9").Value = -42.6813925609733

$ws.Range("A10").Value = "SEDEGI"
$ws.Range("B10").Value = "R. José Maria Figueiredo Réis, 60 - sala 08 - Centro"
$ws.Range("C10").Value = -21.3853296909594
$ws.Range("D10").Value = -42.6903378882941

$ws.Range("A11").Value = "Casa Mattos"
$ws.Range("B11").Value = "R. Ten. Luís Ribeiro, 343 - Vila Domingos Lopes"
$ws.Range("C11").Value = -21.3799180799407
$ws.Range("D11").Value = -42.6906470995146

$ws.Range("A12").Value = "GM2"
$ws.Range("B12").Value = "Av. das Industrias "
$ws.Range("C12").Value = -21.4074040532181
$ws.Range("D12").Value = -42.683095259235

$ws.Range("A13").Value = "Pão Nosso"
$ws.Range("B13").Value = "R. Cel. João Duarte, 120 - Centro"
$ws.Range("C13").Value = -21.3875870108647
$ws.Range("D13").Value = -42.6929887537211

$ws.Range("A14").Value = "Serra Minas - Biscoito Líder de Minas"
$ws.Range("B14").Value = "Avenida Manoel Inacio Peixoto, 300 Distrito Industrial"
$ws.Range("C14").Value = -21.4071550414189
$ws.Range("D14").Value = -42.6808175339905
